$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D-column (price) updates to remain plain text instead of
# being auto-coerced to numbers by COM type inference.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "24.877.47"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").Value = "1.708.68"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "312.48"
$ws.Range("E5").Value = "  +2.14%  "
$ws.Range("D6").Value = "0.9979"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.3760"
$ws.Range("E7").Value = "  +1.20%  "
$ws.Range("D8").Value = "49.65"
$ws.Range("E8").Value = "  +2.82%  "
$ws.Range("D9").Value = "0.3456"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("D11").Value = "0.07568"
$ws.Range("E11").Value = "  +4.02%  "
$ws.Range("D12").Value = "0.9999"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").Value = "21.31"
$ws.Range("E13").Value = "  +4.26%  "
$ws.Range("D14").Value = "6.347"
$ws.Range("E14").Value = "  +3.31%  "
$ws.Range("D15").Value = "7.086"
$ws.Range("E15").Value = "  +4.93%  "
$ws.Range("D16").Value = "1.707.74"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").Value = "0.00001134"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Value = "0.06727"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "0.9984"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "85.33"
$ws.Range("E20").Value = "  +4.84%  "
$ws.Range("D21").Value = "17.41"
$ws.Range("E21").Value = "  +5.59%  "
$ws.Range("E22").Value = "  +4.86%  "
$ws.Range("D23").Value = "13.33"
$ws.Range("E23").Value = "  +11.01%  "
$ws.Range("D24").Value = "24.883.42"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").Value = "2.810"
$ws.Range("E26").Value = "  +4.79%  "
$ws.Range("D27").Value = "20.54"
$ws.Range("E27").Value = "  +4.98%  "
$ws.Range("D28").Value = "151.73"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "133.15"
$ws.Range("E29").Value = "  +4.89%  "
$ws.Range("D30").Value = "1.898.54"
$ws.Range("E30").Value = "  +2.76%  "
$ws.Range("D31").Value = "1.247"
$ws.Range("E31").Value = "  +27.72%  "
$ws.Range("D32").Value = "6.960"
$ws.Range("E32").Value = "  +9.23%  "
$ws.Range("D33").Value = "4.242"
$ws.Range("E33").Value = "  +5.19%  "
$ws.Range("E34").Value = "  +11.74%  "
$ws.Range("D35").Value = "1.794"
$ws.Range("E35").Value = "  +6.47%  "
$ws.Range("D36").Value = "0.08853"
$ws.Range("E36").Value = "  +4.43%  "
$ws.Range("D39").Value = "0.06689"
$ws.Range("E39").Value = "  +2.97%  "
$ws.Range("D40").Value = "0.02416"
$ws.Range("E40").Value = "  +3.27%  "
$ws.Range("D41").Value = "0.2242"
$ws.Range("E41").Value = "  +5.84%  "
$ws.Range("D42").Value = "1.286"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").Value = "0.6471"
$ws.Range("E43").Value = "  +4.42%  "
$ws.Range("D44").Value = "0.9984"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "13.97"
$ws.Range("E45").Value = "  +6.46%  "
$ws.Range("D46").Value = "0.6186"
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("D47").Value = "3.833"
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("D48").Value = "2.142"
$ws.Range("E48").Value = "  +5.26%  "
$ws.Range("D49").Value = "130.55"
$ws.Range("E49").Value = "  +2.67%  "
$ws.Range("D50").Value = "0.07326"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("D51").Value = "80.29"
$ws.Range("E51").Value = "  +5.83%  "

# Rows 37/38: coin order swapped (InternetComputer(DFINITY) now ranks
# above FraxShare) with refreshed price/volume data.
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "5.639"
$ws.Range("E37").Value = "  +5.09%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "9.382"
$ws.Range("E38").Value = "  +4.70%  "

# Restore default style on the price column (clears the temporary
# text-number-format so formatting matches the rest of the sheet).
$priceRange.Style = "Normal"
